# Updates cryptos list: latest price + 1h change figures (GitHub Actions refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '48.079.87'
$ws.Range("E2").Value = '  +1.54%  '
$ws.Range("D3").Value = '2.509.78'
$ws.Range("E3").Value = '  +0.70%  '
$ws.Range("D4").Value = '''0.999'
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").Value = '''321.13'
$ws.Range("E5").Value = '  -0.05%  '
$ws.Range("D6").Value = '''109.13'
$ws.Range("E6").Value = '  +0.91%  '
$ws.Range("E7").Value = '  +1.27%  '
$ws.Range("D8").Value = '''0.999'
$ws.Range("E8").Value = '  -0.07%  '
$ws.Range("D9").Value = '''0.547'
$ws.Range("E9").Value = '  +2.41%  '
$ws.Range("E10").Value = '  +3.09%  '
$ws.Range("D11").Value = '''20.08'
$ws.Range("E11").Value = '  +9.20%  '
$ws.Range("E12").Value = '  +0.96%  '
$ws.Range("D13").Value = '''0.124'
$ws.Range("E13").Value = '  +0.78%  '
$ws.Range("D14").Value = '''7.20'
$ws.Range("E14").Value = '  +1.19%  '
$ws.Range("D15").Value = '2.904.19'
$ws.Range("E15").Value = '  +0.80%  '
$ws.Range("D16").Value = '2.510.59'
$ws.Range("E16").Value = '  +0.74%  '
$ws.Range("D17").Value = '''0.848'
$ws.Range("E17").Value = '  +0.27%  '
$ws.Range("D18").Value = '47.915.28'
$ws.Range("E18").Value = '  +1.40%  '
$ws.Range("D19").Value = '''13.21'
$ws.Range("E20").Value = '  -0.43%  '
$ws.Range("D21").Value = '0.0₃0944'
$ws.Range("E21").Value = '  +1.20%  '
$ws.Range("D22").Value = '''2.71'
$ws.Range("E22").Value = '  +0.37%  '
$ws.Range("D23").Value = '''71.84'
$ws.Range("E23").Value = '  +2.15%  '
$ws.Range("D24").Value = '''276.15'
$ws.Range("E24").Value = '  +12.49%  '
$ws.Range("D25").Value = '''2.56'
$ws.Range("E25").Value = '  -0.57%  '
$ws.Range("E26").Value = '  +0.01%  '
$ws.Range("D27").Value = '''25.90'
$ws.Range("E27").Value = '  +0.59%  '
$ws.Range("D28").Value = '''2.25'
$ws.Range("E28").Value = '  -1.45%  '
$ws.Range("D29").Value = '''10.08'
$ws.Range("E29").Value = '  +0.86%  '
$ws.Range("D30").Value = '''0.142'
$ws.Range("E30").Value = '  +3.40%  '
$ws.Range("D31").Value = '''35.48'
$ws.Range("E31").Value = '  +2.58%  '
$ws.Range("D32").Value = '''49.44'
$ws.Range("E32").Value = '  -0.44%  '
$ws.Range("D33").Value = '''19.44'
$ws.Range("E33").Value = '  -6.61%  '
$ws.Range("D34").Value = '''5.35'
$ws.Range("E34").Value = '  +0.31%  '
$ws.Range("E36").Value = '  -0.09%  '
$ws.Range("D37").Value = '''1.96'
$ws.Range("E37").Value = '  -0.04%  '
$ws.Range("D38").Value = '''4.62'
$ws.Range("E38").Value = '  -0.98%  '
$ws.Range("D39").Value = '''2.96'
$ws.Range("E39").Value = '  +1.35%  '
$ws.Range("B40").Value = 'Stellar'
$ws.Range("C40").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D40").Value = '''0.112'
$ws.Range("E40").Value = '  +0.61%  '
$ws.Range("B41").Value = 'Monero'
$ws.Range("C41").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D41").Value = '''122.38'
$ws.Range("E41").Value = '  +3.85%  '
$ws.Range("E42").Value = '  -0.18%  '
$ws.Range("D43").Value = '''21.83'
$ws.Range("E43").Value = '  -6.06%  '
$ws.Range("D44").Value = '''0.0304'
$ws.Range("E44").Value = '  +2.37%  '
$ws.Range("D45").Value = '2.026.26'
$ws.Range("E45").Value = '  +1.97%  '
$ws.Range("E46").Value = '  +2.30%  '
$ws.Range("E47").Value = '  +4.57%  '
$ws.Range("E48").Value = '  -1.21%  '
$ws.Range("D49").Value = '''9.00'
$ws.Range("E49").Value = '  -0.74%  '
$ws.Range("D50").Value = '''5.16'
$ws.Range("E50").Value = '  +1.15%  '
$ws.Range("D51").Value = '''79.79'
$ws.Range("E51").Value = '  +3.51%  '

foreach ($addr in @('D4', 'D5', 'D6', 'D8', 'D9', 'D11', 'D13', 'D14', 'D17', 'D19', 'D22', 'D23', 'D24', 'D25', 'D27', 'D28', 'D29', 'D30', 'D31', 'D32', 'D33', 'D34', 'D37', 'D38', 'D39', 'D40', 'D41', 'D43', 'D44', 'D49', 'D50', 'D51')) {
    $ws.Range($addr).Style = "Normal"
}
